$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial(-4122)

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 12
$ws.Range("C227").Value = 51
$ws.Range("D227").Value = 154.6251099051026

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 14
$ws.Range("C228").Value = 57
$ws.Range("D228").Value = 172.8162993057029

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 8
$ws.Range("C229").Value = 60
$ws.Range("D229").Value = 181.9118940060031
